$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 22:55"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 211143
$ws.Range("C4").Value = 22613
$ws.Range("E4").Value = 197625
$ws.Range("F4").Value = 5005
$ws.Range("G4").Value = 660
$ws.Range("H4").Value = 4713

# Row 20: Israel -> Brasil
$ws.Range("A20").Value = "Brasil"
$ws.Range("B20").Value = 6836
$ws.Range("C20").Value = 1119
$ws.Range("D20").Value = 127
$ws.Range("E20").Value = 6469
$ws.Range("F20").Value = 296
$ws.Range("G20").Value = 39
$ws.Range("H20").Value = 240

# Row 21: Brasil -> Israel
$ws.Range("A21").Value = "Israel"
$ws.Range("B21").Value = 6092
$ws.Range("C21").Value = 734
$ws.Range("D21").Value = 226
$ws.Range("E21").Value = 5841
$ws.Range("F21").Value = 97
$ws.Range("H21").Value = 25

# Row 32: Polonia -> Polonia
$ws.Range("D32").Value = 47
$ws.Range("E32").Value = 2464

# Row 66: Marruecos -> Marruecos
$ws.Range("B66").Value = 654
$ws.Range("C66").Value = 37
$ws.Range("D66").Value = 29
$ws.Range("E66").Value = 586
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 39

# Row 136: Jamaica -> Jamaica
$ws.Range("E136").Value = 33
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 3

# Row 143: Guam -> El Salvador
$ws.Range("A143").Value = "El Salvador"
$ws.Range("F143").Value = 4

# Row 144: El Salvador -> Guam
$ws.Range("A144").Value = "Guam"
$ws.Range("F144").Value = 0

# Row 149: Tanzania -> Bahamas
$ws.Range("A149").Value = "Bahamas"
$ws.Range("B149").Value = 21
$ws.Range("C149").Value = 7
$ws.Range("E149").Value = 19
$ws.Range("G149").Value = 1

# Row 150: Congo -> Tanzania
$ws.Range("A150").Value = "Tanzania"
$ws.Range("B150").Value = 20
$ws.Range("C150").Value = 1
$ws.Range("D150").Value = 1
$ws.Range("E150").Value = 18
$ws.Range("H150").Value = 1

# Row 151: Maldivas -> Congo
$ws.Range("A151").Value = "Congo"
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 19

# Row 152: Gabon -> Maldivas
$ws.Range("A152").Value = "Maldivas"
$ws.Range("B152").Value = 19
$ws.Range("C152").Value = 1
$ws.Range("D152").Value = 13
$ws.Range("E152").Value = 6
$ws.Range("H152").Value = 0

# Row 153: Islas Virgenes de los Estados Unidos -> Gabon
$ws.Range("A153").Value = "Gabon"
$ws.Range("B153").Value = 18
$ws.Range("C153").Value = 2
$ws.Range("H153").Value = 1

# Row 154: Nueva Caledonia -> Islas Virgenes de los Estados Unidos
$ws.Range("A154").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B154").Value = 17
$ws.Range("E154").Value = 17

# Row 155: Haiti -> Nueva Caledonia
$ws.Range("A155").Value = "Nueva Caledonia"
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 0
$ws.Range("E155").Value = 16

# Row 156: San Martin (Parte Holandesa) -> Haiti
$ws.Range("A156").Value = "Haiti"
$ws.Range("C156").Value = 1
$ws.Range("D156").Value = 1
$ws.Range("E156").Value = 15
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 0

# Row 157: Eritrea -> San Martin (Parte Holandesa)
$ws.Range("A157").Value = "San Martin (Parte Holandesa)"
$ws.Range("B157").Value = 16
$ws.Range("C157").Value = 10
$ws.Range("D157").Value = 6
$ws.Range("E157").Value = 9
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 1

# Row 158: Birmania -> Eritrea
$ws.Range("A158").Value = "Eritrea"
$ws.Range("E158").Value = 15
$ws.Range("H158").Value = 0

# Row 160: Bahamas -> Birmania
$ws.Range("A160").Value = "Birmania"
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 0
$ws.Range("H160").Value = 1

# Row 163: Mongolia -> Namibia
$ws.Range("A163").Value = "Namibia"
$ws.Range("C163").Value = 3

# Row 164: Namibia -> Mongolia
$ws.Range("A164").Value = "Mongolia"
$ws.Range("C164").Value = 2
